# Swipe to toggle done in TodoListActivity (62/100)
#
# - The "Erledigtsein/Nicht-Erledigtsein" toggle item (row 32) is now fully
#   implemented (swipe-to-toggle), so it gets full credit (D32: 1 -> 3) and
#   the reviewer comment that it didn't seem to work yet is removed (E32).
# - The "Todoliste" section marker (column A) actually belongs on row 25
#   (where the "Todoliste" overview description starts), not row 26, so it
#   is moved up one row.
# - Row 24's comment text got shorter, so its custom row height shrinks.
# - The overall score total (D55) recalculates automatically from 60 to 62.
# - Selection / scroll position is updated to reflect where the author was
#   last working (around D32).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24: comment box is shorter now -> smaller custom row height.
$ws.Rows.Item(24).RowHeight = 104.4

# Move the "Todoliste" category marker from A26 up to A25 (cut/paste keeps
# both the shared-string value and the cell style intact).
$ws.Range("A26").Cut($ws.Range("A25")) | Out-Null
$ws.Range("A26").Clear() | Out-Null

# Row 32: "Erledigtsein/Nicht-Erledigtsein" toggle now works -> full marks,
# and drop the now-stale "doesn't seem to work yet" reviewer comment.
$ws.Range("D32").Value = 3
$ws.Range("E32").Clear() | Out-Null

# Reflect the last-looked-at cell / scroll position.
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D32").Select() | Out-Null
